$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 12:52"

# Etiopia overtakes Madagascar in ranking (row 141 becomes Etiopia, row 142 becomes Madagascar)
$ws.Range("A141").Value = "Etiopia"
$ws.Range("A142").Value = "Madagascar"

# Row 9 - Alemania
$ws.Range("B9").Value = 160059
$ws.Range("C9").Value = 147
$ws.Range("E9").Value = 33345

# Row 12 - Iran
$ws.Range("B12").Value = 93657
$ws.Range("C12").Value = 1073
$ws.Range("D12").Value = 73791
$ws.Range("E12").Value = 13909
$ws.Range("F12").Value = 2965
$ws.Range("G12").Value = 80
$ws.Range("H12").Value = 5957

# Row 20 - Suiza
$ws.Range("B20").Value = 29407
$ws.Range("C20").Value = 143
$ws.Range("E20").Value = 5108

# Row 47 - Australia
$ws.Range("B47").Value = 6746
$ws.Range("C47").Value = 8
$ws.Range("D47").Value = 5667
$ws.Range("E47").Value = 990
$ws.Range("F47").Value = 38
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 89

# Row 68 - Uzbekistan
$ws.Range("B68").Value = 1969
$ws.Range("C68").Value = 30
$ws.Range("E68").Value = 906

# Row 99 - Libano
$ws.Range("B99").Value = 721
$ws.Range("C99").Value = 4
$ws.Range("D99").Value = 150
$ws.Range("E99").Value = 547

# Row 129 - Maldivas
$ws.Range("B129").Value = 256
$ws.Range("C129").Value = 6
$ws.Range("E129").Value = 239

# Row 141 - Etiopia (updated stats, new rank above Madagascar)
$ws.Range("B141").Value = 130
$ws.Range("C141").Value = 4
$ws.Range("D141").Value = 58
$ws.Range("E141").Value = 69
$ws.Range("F141").Value = 0
$ws.Range("H141").Value = 3

# Row 142 - Madagascar (stats unchanged, moved down one rank)
$ws.Range("B142").Value = 128
$ws.Range("D142").Value = 82
$ws.Range("E142").Value = 46
$ws.Range("F142").Value = 1
$ws.Range("H142").Value = 0

# Row 167 - Nepal
$ws.Range("B167").Value = 57
$ws.Range("C167").Value = 3
$ws.Range("E167").Value = 41
